# Sprint 11 alignment: add procedureName / procedureDescription columns to
# the studyDesignProcedures sheet, populate two new test rows, and make
# that sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesignProcedures")

# Insert two new blank columns before the existing "procedureCode" column
# (currently column C), pushing procedureCode / procedureIsConditional /
# procedureIsConditionalReason two columns to the right.
$ws.Range("C1:D1").EntireColumn.Insert()

# Match the column widths the author ended up with (character units, which
# this host stores with a small fixed padding offset vs the raw ColumnWidth
# value, hence the slightly odd-looking constants below).
$ws.Range("C1").ColumnWidth = 16.666666666666668
$ws.Range("D1").ColumnWidth = 22.998697916666668

# Header cells: use a leading apostrophe so the text keeps the same
# "quote prefix" flavoured style the rest of the header row already uses.
$ws.Range("C1").Formula = "'procedureName"
$ws.Range("D1").Formula = "'procedureDescription"

# New data rows (plain values -- no quote-prefix styling on these).
$ws.Range("C2").Value = "Test8"
$ws.Range("D2").Value = "Test Eight"
$ws.Range("C3").Value = "Test9"
$ws.Range("D3").Value = "Test Nine"

# Reflect the author's final selection and make this sheet the active tab.
$ws.Range("C1:D3").Select() | Out-Null
$ws.Activate() | Out-Null
